# Generate Report for Handoff
# Updates the localization-status report for the "b.md" file: it is now
# ready for handoff, with a freshly generated handoff xliff and a note
# that the previous handback was based on a stale source version.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5ad9491e589e8f3173262c0e91f33a1557016580/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6130edbeb2d0dd5b8802f44c6e7172fa4efe655c/e2e/b.md."

# ---- Overview sheet ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-09-04 08:41:08"

# ---- zh-cn sheet ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Columns.Item(16).ColumnWidth = 39.16
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zh.Range("H3").Value = "2016-09-04 08:41:00"
$zh.Range("P3").Value = $errorDetail

# ---- de-de sheet ----
$de = $wb.Worksheets.Item("de-de")
$de.Columns.Item(16).ColumnWidth = 39.16
$de.Range("C3").Value = "Ready for handoff"
$de.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$de.Range("H3").Value = "2016-09-04 08:41:08"
$de.Range("P3").Value = $errorDetail

Write-Output "Applied handoff report updates"
